# Daily attendance processing - 2026-01-30 21:12:58
# Rotate the "Recorded By" (column G) comma-separated list of recorders
# left by one position (move the first entry to the end) for every row
# that has more than one recorder listed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $v = $cell.Value2

    if ($v -ne $null -and $v -like "*,*") {
        $parts = $v -split ", "
        $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
        $cell.Value2 = $rotated -join ", "
    }
}
